# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows whose data changed on repull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -5
    7  = -2
    13 = -3
    14 = -1
    15 = -1
    16 = -1
    17 = 1
    23 = -1
    25 = 0
    28 = -2
    38 = 0
    40 = 3
    41 = 1
    44 = 1
    45 = 1
    56 = -5
    57 = 0
    58 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
